{"js": "// Edit for \"Examples for ALNfitDeep.docx\"\n// Rewrites the tail of paragraph I, splits off a new \"Advanced:\" paragraph,\n// rewrites paragraph J, and adds two new trailing paragraphs (\"...DTREEs\n// are not smoothed.\" and \"Many other improvements...\").\n//\n// The paragraph-J edit is done in separate search/replace passes split at\n// the pre-existing \"_GoBack\" bookmark's original location so that bookmark\n// keeps its place in the surrounding text (it ends up between \"M\" and\n// \"any other \", matching the target document) instead of being swallowed\n// by a single large replace.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------\n// 1) Paragraph I: replace the tail text (from \"You can try...\" through\n//    the old \"(advanced: ... replacements).\" sentence) with the new\n//    wording that ends with \"...before training.  \"\n// ---------------------------------------------------------------\nconst oldITail =\n  \" You can try more complicated real-world data where there are missing \" +\n  \"values in several columns.  The missing values may be replaced by \" +\n  \"several trainings where the values in columns not removed are \" +\n  \"completely defined.  (advanced: Use knowledge of the functional \" +\n  \"dependencies in a relational database to determine the sequence of \" +\n  \"column replacements).\";\n\nconst newITail =\n  \" You can try more complicated real-world data where there are missing \" +\n  \"values in several columns.  The missing values may be replaced by \" +\n  \"successive trainings where the values in input columns are completely \" +\n  \"defined.  Just select the last R output file as the next datafile in a \" +\n  \"run.  Alternatively, inputs in a column can be defined if you remove \" +\n  \"the rows with missing values before training.  \";\n\nlet results = body.search(oldITail, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\"Paragraph I target text not found (count=\" + results.items.length + \")\");\n}\nlet target = results.items[0];\ntarget.insertText(newITail, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------\n// 2) Insert the new \"Advanced:\" paragraph right after paragraph I.\n// ---------------------------------------------------------------\nconst advancedText =\n  \"Advanced:  If the file with missing values is a relation in a \" +\n  \"relational database, you can use knowledge of the functional \" +\n  \"dependencies in that relation to suggest a sequence of column \" +\n  \"replacements.\";\n\nresults = body.search(\"before training.  \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\"End-of-paragraph-I anchor not found (count=\" + results.items.length + \")\");\n}\ntarget = results.items[0];\ntarget.insertParagraph(advancedText, Word.InsertLocation.after);\nawait context.sync();\n\n// ---------------------------------------------------------------\n// 3) Paragraph J, part 1: replace everything up to (but not including)\n//    the \"improvements...\" text that follows the \"_GoBack\" bookmark.\n//    Embedded \"\\n\" characters create the two new paragraph breaks that\n//    turn the old paragraph J into three paragraphs: the rewritten J, a\n//    new paragraph about ALNfitDeep/DTREEs, and the start (\"M\") of the\n//    final paragraph.\n// ---------------------------------------------------------------\nconst oldJPre =\n  \"J. Deep learning is faster and allows more complicated functions to \" +\n  \"be learned.  Smoothing allows one to use fewer flat pieces in the \" +\n  \"learned function, making it easier to analyze.  In this software, the \" +\n  \"only step to promote deep learning is that smoothing can be switched \" +\n  \"off, however many other \";\n\nconst newJPre =\n  \"J. .  Complicated functions will be approximated by many flat pieces \" +\n  \"in the DTREE produced in each run.  Smoothing allows one to use fewer \" +\n  \"flat pieces in the learned function for the same accuracy, making the \" +\n  \"learned result easier to analyze by analyzing the flat pieces.  The \" +\n  \"deviation of a smoothed ALN from its ALN without smoothing depends on \" +\n  \"a constant, the Smoothing Epsilon, and the depth of the ALN in terms \" +\n  \"of maximum or minimum nodes.  Some small changes are required in \" +\n  \"ALNfitDeep to be able to make sure the learned function satisfies a \" +\n  \"given specification based on the analysis of the flat pieces. \\n\" +\n  \"In ALNfitDeep, the only step to promote deep learning is that \" +\n  \"smoothing can be set to zero.  Deep learning is faster and allows \" +\n  \"more complicated functions to be learned by deeper ALNs  DTREEs are \" +\n  \"not smoothed.\\n\";\n\nresults = body.search(oldJPre, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\"Paragraph J (pre-bookmark) text not found (count=\" + results.items.length + \")\");\n}\ntarget = results.items[0];\ntarget.insertText(newJPre, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------\n// 4) Insert the \"M\" that begins the final paragraph right before the\n//    (still present) \"_GoBack\" bookmark / \"improvements...\" text.\n// ---------------------------------------------------------------\nconst oldJPost =\n  \"improvements could be made using this model of deep learning. This \" +\n  \"is just an indication of what is possible.\";\n\nresults = body.search(oldJPost, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\"Paragraph J (post-bookmark) text not found, pass 1 (count=\" + results.items.length + \")\");\n}\ntarget = results.items[0];\ntarget.insertText(\"M\", Word.InsertLocation.before);\nawait context.sync();\n\n// ---------------------------------------------------------------\n// 5) Replace the text after the bookmark with its final wording.\n// ---------------------------------------------------------------\nconst newJPost =\n  \"any other improvements could be made using this model of ALN deep \" +\n  \"learning. The present software is just an indication of what is \" +\n  \"possible.\";\n\nresults = body.search(oldJPost, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length !== 1) {\n  throw new Error(\"Paragraph J (post-bookmark) text not found, pass 2 (count=\" + results.items.length + \")\");\n}\ntarget = results.items[0];\ntarget.insertText(newJPost, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Edit for \"Examples for ALNfitDeep.docx\"\n# Rewrites the tail of paragraph I, splits off a new \"Advanced:\" paragraph,\n# rewrites paragraph J, and adds two new trailing paragraphs (\"...DTREEs are\n# not smoothed.\" and \"Many other improvements...\").\n#\n# Uses Range.Find.Execute with \"^p\" in the replacement text to insert new\n# paragraph marks at the split points. The paragraph-J edit is done in two\n# Find/Replace passes (before / after the \"_GoBack\" bookmark) plus a small\n# InsertBefore so the pre-existing bookmark keeps its place in the\n# surrounding text instead of being swallowed by a Find/Replace.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------\n# 1) Paragraph I: replace the tail text (from \"You can try...\" through\n#    the old \"(advanced: ... replacements).\" sentence) with the new\n#    wording, then append the new \"Advanced:\" paragraph right after it.\n# ---------------------------------------------------------------\n$oldITail = \" You can try more complicated real-world data where there are missing values in several columns.  The missing values may be replaced by several trainings where the values in columns not removed are completely defined.  (advanced: Use knowledge of the functional dependencies in a relational database to determine the sequence of column replacements).\"\n$newITail = \" You can try more complicated real-world data where there are missing values in several columns.  The missing values may be replaced by successive trainings where the values in input columns are completely defined.  Just select the last R output file as the next datafile in a run.  Alternatively, inputs in a column can be defined if you remove the rows with missing values before training.  ^pAdvanced:  If the file with missing values is a relation in a relational database, you can use knowledge of the functional dependencies in that relation to suggest a sequence of column replacements.\"\n\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$found1 = $rng1.Find.Execute($oldITail, $false, $false, $false, $false, $false, $true, 1, $false, $newITail, 2)\nWrite-Output $found1\n\n# ---------------------------------------------------------------\n# 2) Paragraph J: replace everything up to (but not including) the\n#    \"improvements...\" text that follows the \"_GoBack\" bookmark. This\n#    also inserts two paragraph breaks, turning the tail of the old\n#    paragraph J into three paragraphs: the rewritten J, a new paragraph\n#    about ALNfitDeep/DTREEs, and the start (\"M\") of the final paragraph.\n# ---------------------------------------------------------------\n$oldJPreBookmark = \"J. Deep learning is faster and allows more complicated functions to be learned.  Smoothing allows one to use fewer flat pieces in the learned function, making it easier to analyze.  In this software, the only step to promote deep learning is that smoothing can be switched off, however many other \"\n$newJPreBookmark = \"J. .  Complicated functions will be approximated by many flat pieces in the DTREE produced in each run.  Smoothing allows one to use fewer flat pieces in the learned function for the same accuracy, making the learned result easier to analyze by analyzing the flat pieces.  The deviation of a smoothed ALN from its ALN without smoothing depends on a constant, the Smoothing Epsilon, and the depth of the ALN in terms of maximum or minimum nodes.  Some small changes are required in ALNfitDeep to be able to make sure the learned function satisfies a given specification based on the analysis of the flat pieces. ^pIn ALNfitDeep, the only step to promote deep learning is that smoothing can be set to zero.  Deep learning is faster and allows more complicated functions to be learned by deeper ALNs  DTREEs are not smoothed.^p\"\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$found2 = $rng2.Find.Execute($oldJPreBookmark, $false, $false, $false, $false, $false, $true, 1, $false, $newJPreBookmark, 2)\nWrite-Output $found2\n\n# ---------------------------------------------------------------\n# 2b) Insert the \"M\" that begins the final paragraph right before the\n#     (still present) \"_GoBack\" bookmark / \"improvements...\" text, so the\n#     bookmark ends up sitting between \"M\" and \"any other \" as in the\n#     target document.\n# ---------------------------------------------------------------\n$oldJPostBookmark = \"improvements could be made using this model of deep learning. This is just an indication of what is possible.\"\n\n$rng2b = $d.Content\n$rng2b.Find.ClearFormatting()\n$found2b = $rng2b.Find.Execute($oldJPostBookmark)\nWrite-Output $found2b\n$rng2b.Collapse(1)  # wdCollapseStart\n$rng2b.InsertBefore(\"M\")\n\n# ---------------------------------------------------------------\n# 3) Replace the text after the bookmark with its final wording.\n# ---------------------------------------------------------------\n$newJPostBookmark = \"any other improvements could be made using this model of ALN deep learning. The present software is just an indication of what is possible.\"\n\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Replacement.ClearFormatting()\n$found3 = $rng3.Find.Execute($oldJPostBookmark, $false, $false, $false, $false, $false, $true, 1, $false, $newJPostBookmark, 2)\nWrite-Output $found3\n"}
